$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.449.83'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '1.698.94'
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.02'
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5480'
$ws.Range("E6").Value = '  +3.96%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2741'
$ws.Range("E8").Value = '  +1.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06448'
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.07'
$ws.Range("E11").Value = '  +3.00%  '
$ws.Range("D12").Value = '1.750.58'
$ws.Range("E12").Value = '  +3.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.553'
$ws.Range("E13").Value = '  +0.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5846'
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008431'
$ws.Range("E15").Value = '  -0.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.81'
$ws.Range("E16").Value = '  +2.33%  '
$ws.Range("D17").Value = '26.510.76'
$ws.Range("E17").Value = '  +0.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.952'
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.00'
$ws.Range("E20").Value = '  +1.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.62'
$ws.Range("E21").Value = '  +1.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.265'
$ws.Range("E22").Value = '  +1.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.010'
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("E24").Value = '  +3.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1324'
$ws.Range("E25").Value = '  +6.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.903'
$ws.Range("E26").Value = '  +2.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.81'
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06266'
$ws.Range("E28").Value = '  -5.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.380'
$ws.Range("E29").Value = '  +2.53%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.616'
$ws.Range("E31").Value = '  +1.50%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.601'
$ws.Range("E32").Value = '  +0.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.697'
$ws.Range("E33").Value = '  +2.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.042'
$ws.Range("E34").Value = '  +1.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6184'
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.412'
$ws.Range("E36").Value = '  +0.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.770'
$ws.Range("E37").Value = '  +2.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01647'
$ws.Range("E38").Value = '  +1.48%  '
$ws.Range("D39").Value = '1.119.42'
$ws.Range("E39").Value = '  +0.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.143'
$ws.Range("E40").Value = '  -3.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8776'
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("E43").Value = '  +0.76%  '
$ws.Range("D44").Value = '1.851.28'
$ws.Range("E44").Value = '  +1.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000111'
$ws.Range("E45").Value = '  -3.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.63'
$ws.Range("E46").Value = '  +1.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.253'
$ws.Range("E47").Value = '  +0.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.006'
$ws.Range("E48").Value = '  +0.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05284'
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.155'
$ws.Range("E50").Value = '  +1.79%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4302'
$ws.Range("E51").Value = '  -0.09%  '
